$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "62.688.82"
$ws.Range("E2").Value = "  -1.87%  "

# Row 3
$ws.Range("D3").Value = "2.443.15"
$ws.Range("E3").Value = "  -1.79%  "

# Row 4
$ws.Range("D4").Formula = "'0.999"
$ws.Range("E4").Value = "  -0.32%  "

# Row 5
$ws.Range("D5").Formula = "'571.25"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("D6").Formula = "'145.53"
$ws.Range("E6").Value = "  -3.37%  "

# Row 7
$ws.Range("E7").Value = "  +0.22%  "

# Row 8
$ws.Range("D8").Formula = "'0.533"
$ws.Range("E8").Value = "  -1.39%  "

# Row 9
$ws.Range("D9").Value = "2.439.80"
$ws.Range("E9").Value = "  -2.49%  "

# Row 10
$ws.Range("D10").Formula = "'0.110"
$ws.Range("E10").Value = "  -4.51%  "

# Row 11
$ws.Range("D11").Formula = "'0.157"
$ws.Range("E11").Value = "  +1.52%  "

# Row 12
$ws.Range("D12").Formula = "'5.25"
$ws.Range("E12").Value = "  -2.35%  "

# Row 13
$ws.Range("D13").Formula = "'0.356"
$ws.Range("E13").Value = "  -3.21%  "

# Row 14
$ws.Range("D14").Formula = "'27.06"
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").Formula = "'0.0000177"
$ws.Range("E15").Value = "  -5.62%  "

# Row 16
$ws.Range("D16").Value = "2.884.72"
$ws.Range("E16").Value = "  -1.75%  "

# Row 17
$ws.Range("D17").Value = "62.565.38"
$ws.Range("E17").Value = "  -1.36%  "

# Row 18
$ws.Range("D18").Value = "2.455.01"
$ws.Range("E18").Value = "  -2.20%  "

# Row 19
$ws.Range("D19").Formula = "'11.28"
$ws.Range("E19").Value = "  -2.59%  "

# Row 20
$ws.Range("E20").Value = "  -2.09%  "

# Row 21
$ws.Range("D21").Formula = "'327.82"
$ws.Range("E21").Value = "  -0.68%  "

# Row 22
$ws.Range("D22").Formula = "'4.17"
$ws.Range("E22").Value = "  -1.84%  "

# Row 23
$ws.Range("D23").Formula = "'2.11"
$ws.Range("E23").Value = "  +10.23%  "

# Row 24
$ws.Range("D24").Formula = "'1.00"
$ws.Range("E24").Value = "  +0.09%  "

# Row 25
$ws.Range("D25").Formula = "'65.32"
$ws.Range("E25").Value = "  -3.72%  "

# Row 26
$ws.Range("D26").Formula = "'623.21"
$ws.Range("E26").Value = "  -5.52%  "

# Row 27
$ws.Range("D27").Formula = "'8.98"
$ws.Range("E27").Value = "  -0.07%  "

# Row 28
$ws.Range("D28").Formula = "'0.0000100"
$ws.Range("E28").Value = "  -6.22%  "

# Row 29
$ws.Range("D29").Value = "2.561.62"
$ws.Range("E29").Value = "  -1.52%  "

# Row 30
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Formula = "'1.50"
$ws.Range("E30").Value = "  -3.13%  "

# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Formula = "'1.00"
$ws.Range("E31").Value = "  +0.53%  "

# Row 32
$ws.Range("D32").Formula = "'8.15"
$ws.Range("E32").Value = "  -5.70%  "

# Row 33
$ws.Range("D33").Formula = "'1.89"
$ws.Range("E33").Value = "  -2.39%  "

# Row 34
$ws.Range("D34").Formula = "'0.139"
$ws.Range("E34").Value = "  -5.32%  "

# Row 35
$ws.Range("D35").Formula = "'5.15"
$ws.Range("E35").Value = "  -1.90%  "

# Row 36
$ws.Range("D36").Formula = "'1.50"
$ws.Range("E36").Value = "  -4.48%  "

# Row 37
$ws.Range("E37").Value = "  +0.27%  "

# Row 38
$ws.Range("D38").Formula = "'0.378"
$ws.Range("E38").Value = "  -3.10%  "

# Row 39
$ws.Range("D39").Formula = "'18.80"
$ws.Range("E39").Value = "  -1.69%  "

# Row 40
$ws.Range("D40").Formula = "'5.33"
$ws.Range("E40").Value = "  -5.72%  "

# Row 41
$ws.Range("D41").Formula = "'146.40"
$ws.Range("E41").Value = "  -1.55%  "

# Row 42
$ws.Range("D42").Formula = "'1.78"
$ws.Range("E42").Value = "  -5.05%  "

# Row 43
$ws.Range("D43").Formula = "'2.60"
$ws.Range("E43").Value = "  -3.37%  "

# Row 44
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Formula = "'146.83"
$ws.Range("E45").Value = "  -4.77%  "

# Row 46
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Formula = "'3.77"
$ws.Range("E46").Value = "  -0.61%  "

# Row 47
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Formula = "'20.82"
$ws.Range("E47").Value = "  -2.82%  "

# Row 48
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Formula = "'0.0529"
$ws.Range("E48").Value = "  -4.32%  "

# Row 49
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Formula = "'0.597"
$ws.Range("E49").Value = "  -2.96%  "

# Row 50
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Formula = "'0.0232"
$ws.Range("E50").Value = "  -3.32%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Formula = "'0.0918"
$ws.Range("E51").Value = "  -1.39%  "
